$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 12416.667
$ws.Range("J7").Value = 12416.667
$ws.Range("L7").Value = 12416.667
$ws.Range("N7").Value = -12640.667

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 12416.667
$ws.Range("J14").Value = 12416.667
$ws.Range("L14").Value = 12416.667
$ws.Range("N14").Value = -12798.667

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 704.7037
$ws.Range("J28").Value = 1089.7142
$ws.Range("L28").Value = 1089.7142
$ws.Range("N28").Value = -2059.7142

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2080.4614
$ws.Range("I33").Value = 607.2857
$ws.Range("K33").Value = 607.2857
$ws.Range("M33").Value = -378.2857

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3367.8572
$ws.Range("I141").Value = 2771.25
$ws.Range("J141").Value = 4163.3335
$ws.Range("K141").Value = 8313.75
$ws.Range("L141").Value = 12490.0005
$ws.Range("M141").Value = -3133.75
$ws.Range("N141").Value = -22850.0005

# ARM row 6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 6818.4546
$ws.Range("I6").Value = 5000.6665
$ws.Range("J6").Value = 8999.799999999999
$ws.Range("K6").Value = 5000.6665
$ws.Range("L6").Value = 8999.799999999999
$ws.Range("M6").Value = -4827.6665
$ws.Range("N6").Value = -9345.799999999999

# ARM row 25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3283.2
$ws.Range("I25").Value = 708
$ws.Range("K25").Value = 708
$ws.Range("M25").Value = -306

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1875.1875
$ws.Range("J88").Value = 2137.4
$ws.Range("L88").Value = 2137.4
$ws.Range("N88").Value = -2949.4

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1875.1875
$ws.Range("J91").Value = 2137.4
$ws.Range("L91").Value = 2137.4
$ws.Range("N91").Value = -4945.4

# BSM row 25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 49375.176
$ws.Range("I86").Value = 74638
$ws.Range("K86").Value = 74638
$ws.Range("M86").Value = -73515

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 49375.176
$ws.Range("I89").Value = 74638
$ws.Range("K89").Value = 373190
$ws.Range("M89").Value = -367574

# CRP row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 20
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 20
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("M3").Value = 93

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1007.36365
$ws.Range("I16").Value = 820
$ws.Range("J16").Value = 1163.5
$ws.Range("K16").Value = 820
$ws.Range("L16").Value = 1163.5
$ws.Range("M16").Value = -533
$ws.Range("N16").Value = -1737.5

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2598.6365
$ws.Range("I62").Value = 2498
$ws.Range("K62").Value = 2498
$ws.Range("M62").Value = -1874

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2598.6365
$ws.Range("I65").Value = 2498
$ws.Range("K65").Value = 12490
$ws.Range("M65").Value = -9370

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1007.36365
$ws.Range("I113").Value = 820
$ws.Range("J113").Value = 1163.5
$ws.Range("K113").Value = 820
$ws.Range("L113").Value = 1163.5
$ws.Range("M113").Value = 1350
$ws.Range("N113").Value = -5503.5

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 69.2
$ws.Range("I6").Value = 69.2
$ws.Range("K6").Value = 207.6
$ws.Range("M6").Value = -94.60000000000002

# CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2000
$ws.Range("J64").Value = 2000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6540

# CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2000
$ws.Range("J67").Value = 2000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7872

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1970.2727
$ws.Range("I107").Value = 1397.6666
$ws.Range("J107").Value = 2185
$ws.Range("K107").Value = 4192.9998
$ws.Range("L107").Value = 6555
$ws.Range("M107").Value = -2272.9998
$ws.Range("N107").Value = -10395

# CUL row 110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 6300
$ws.Range("I110").Value = 2500
$ws.Range("K110").Value = 7500
$ws.Range("M110").Value = -3410

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 776.7778
$ws.Range("I131").Value = 355.2353
$ws.Range("J131").Value = 864.1707
$ws.Range("K131").Value = 1065.7059
$ws.Range("L131").Value = 2592.5121
$ws.Range("M131").Value = 3974.2941
$ws.Range("N131").Value = -12672.5121

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2717.6924
$ws.Range("I141").Value = 1842.2222
$ws.Range("J141").Value = 4687.5
$ws.Range("K141").Value = 5526.6666
$ws.Range("L141").Value = 14062.5
$ws.Range("M141").Value = -346.6665999999996
$ws.Range("N141").Value = -24422.5

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 79071.55499999999
$ws.Range("I70").Value = 171221.5
$ws.Range("J70").Value = 5351.6
$ws.Range("K70").Value = 171221.5
$ws.Range("L70").Value = 5351.6
$ws.Range("M70").Value = -170951.5
$ws.Range("N70").Value = -5891.6

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 79071.55499999999
$ws.Range("I73").Value = 171221.5
$ws.Range("J73").Value = 5351.6
$ws.Range("K73").Value = 171221.5
$ws.Range("L73").Value = 5351.6
$ws.Range("M73").Value = -170285.5
$ws.Range("N73").Value = -7223.6

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2264852
$ws.Range("I126").Value = 2789
$ws.Range("J126").Value = 3678641.2
$ws.Range("K126").Value = 8367
$ws.Range("L126").Value = 11035923.6
$ws.Range("M126").Value = -5897
$ws.Range("N126").Value = -11040863.6

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2590
$ws.Range("I7").Value = 1975
$ws.Range("J7").Value = 3820
$ws.Range("K7").Value = 1975
$ws.Range("L7").Value = 3820
$ws.Range("M7").Value = -1863
$ws.Range("N7").Value = -4044

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2590
$ws.Range("I126").Value = 1975
$ws.Range("J126").Value = 3820
$ws.Range("K126").Value = 5925
$ws.Range("L126").Value = 11460
$ws.Range("M126").Value = -3455
$ws.Range("N126").Value = -16400

# WVR row 26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 169000
$ws.Range("J26").Value = 252500
$ws.Range("L26").Value = 252500
$ws.Range("N26").Value = -253086

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 4385.9443
$ws.Range("I54").Value = 7070
$ws.Range("J54").Value = 4228.0586
$ws.Range("K54").Value = 7070
$ws.Range("L54").Value = 4228.0586
$ws.Range("M54").Value = -6550
$ws.Range("N54").Value = -5268.0586

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 286865.84
$ws.Range("J81").Value = 168010.17
$ws.Range("L81").Value = 336020.34
$ws.Range("N81").Value = -338142.34

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 286865.84
$ws.Range("J84").Value = 168010.17
$ws.Range("L84").Value = 1680101.7
$ws.Range("N84").Value = -1690709.7
